$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns stay formatted as text so values like
# "1.00" or "0.0340" are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "62.726.00"
$ws.Range("E2").Value = "  +2.66%  "

# Row 3
$ws.Range("D3").Value = "2.967.36"
$ws.Range("E3").Value = "  +1.33%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "594.31"
$ws.Range("E5").Value = "  +0.26%  "

# Row 6
$ws.Range("D6").Value = "145.82"
$ws.Range("E6").Value = "  +0.54%  "

# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").Value = "2.967.70"
$ws.Range("E8").Value = "  +1.37%  "

# Row 9
$ws.Range("E9").Value = "  +0.40%  "

# Row 10
$ws.Range("D10").Value = "7.23"
$ws.Range("E10").Value = "  +2.90%  "

# Row 11
$ws.Range("D11").Value = "0.144"
$ws.Range("E11").Value = "  +0.90%  "

# Row 12
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  +0.67%  "

# Row 13
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  +5.42%  "

# Row 14
$ws.Range("D14").Value = "33.29"
$ws.Range("E14").Value = "  -1.17%  "

# Row 15
$ws.Range("E15").Value = "  -0.33%  "

# Row 16
$ws.Range("D16").Value = "3.457.50"
$ws.Range("E16").Value = "  +1.28%  "

# Row 17
$ws.Range("D17").Value = "62.638.27"
$ws.Range("E17").Value = "  +2.51%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.005.33"
$ws.Range("E18").Value = "  +2.64%  "

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "6.70"
$ws.Range("E19").Value = "  -0.57%  "

# Row 20
$ws.Range("D20").Value = "442.12"
$ws.Range("E20").Value = "  +1.64%  "

# Row 21
$ws.Range("D21").Value = "13.43"
$ws.Range("E21").Value = "  -0.22%  "

# Row 22
$ws.Range("D22").Value = "0.672"
$ws.Range("E22").Value = "  -1.03%  "

# Row 23
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("D24").Value = "81.88"
$ws.Range("E24").Value = "  +0.42%  "

# Row 25
$ws.Range("E25").Value = "  +2.15%  "

# Row 26
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "2.14"
$ws.Range("E26").Value = "  -3.00%  "

# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "11.87"
$ws.Range("E27").Value = "  +0.27%  "

# Row 28
$ws.Range("E28").Value = "  -0.01%  "

# Row 29
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  +3.60%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.60"
$ws.Range("E30").Value = "  -0.20%  "

# Row 31
$ws.Range("E31").Value = "  -5.49%  "

# Row 32
$ws.Range("D32").Value = "26.67"
$ws.Range("E32").Value = "  -0.19%  "

# Row 33
$ws.Range("D33").Value = "0.110"
$ws.Range("E33").Value = "  -1.11%  "

# Row 34
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0930"
$ws.Range("E34").Value = "  +7.35%  "

# Row 35
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -1.72%  "

# Row 37
$ws.Range("D37").Value = "5.63"
$ws.Range("E37").Value = "  -0.32%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "2.99"
$ws.Range("E38").Value = "  +0.22%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "2.04"
$ws.Range("E39").Value = "  +2.52%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "49.53"
$ws.Range("E40").Value = "  -0.86%  "

# Row 41
$ws.Range("D41").Value = "8.56"
$ws.Range("E41").Value = "  -0.28%  "

# Row 42
$ws.Range("D42").Value = "0.119"
$ws.Range("E42").Value = "  -4.69%  "

# Row 43
$ws.Range("D43").Value = "0.282"
$ws.Range("E43").Value = "  -0.97%  "

# Row 44
$ws.Range("D44").Value = "39.35"
$ws.Range("E44").Value = "  -6.78%  "

# Row 45
$ws.Range("D45").Value = "2.738.83"
$ws.Range("E45").Value = "  +1.78%  "

# Row 46
$ws.Range("D46").Value = "136.17"
$ws.Range("E46").Value = "  +1.63%  "

# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "365.54"
$ws.Range("E47").Value = "  -2.27%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0340"
$ws.Range("E48").Value = "  -1.92%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "23.08"
$ws.Range("E50").Value = "  -3.19%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").Value = "  -0.32%  "
